# Apply bug-fix updates to the "Impressions"-related figures on rows 13 and 15
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Impressions (D13) 5020 -> 5021
$ws.Range("D13").Value = 5021

# Row 15: Impressions (D15) 704 -> 708
$ws.Range("D15").Value = 708

# Row 15: Impressions Organic (F15) 290 -> 293
$ws.Range("F15").Value = 293

# Row 15: Impressions Viral (J15) 184 -> 188
$ws.Range("J15").Value = 188
